$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear old totals row (row 13) before rebuilding the sheet with new data rows
$ws.Range("B13:C13").ClearContents()

# New data rows 6-8 (records 5-7): still being chased up, so the "date" column
# just holds a placeholder ellipsis string, styled like the rest of column D
# (centered, general number format -> same cellXf as the column default).
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = 10000
$ws.Cells.Item(6,3).Value = "Thang Duong"
$ws.Cells.Item(6,4).HorizontalAlignment = -4108
$ws.Cells.Item(6,4).Value = "….."
$ws.Cells.Item(6,5).Value = "late for meeting"
$ws.Cells.Item(6,6).Value = "Waiting"

$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = 10000
$ws.Cells.Item(7,3).Value = "Thi Nguyen"
$ws.Cells.Item(7,4).HorizontalAlignment = -4108
$ws.Cells.Item(7,4).Value = "…"
$ws.Cells.Item(7,5).Value = "late for meeting"
$ws.Cells.Item(7,6).Value = "Waiting"

$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = 10000
$ws.Cells.Item(8,3).Value = "Thang Duong"
$ws.Cells.Item(8,4).HorizontalAlignment = -4108
$ws.Cells.Item(8,4).Value = "…"
$ws.Cells.Item(8,5).Value = "late for meeting"
$ws.Cells.Item(8,6).Value = "Waiting"

# Rows 9-14 (records 8-13): real dates, formatted like the existing D2:D5 cells
# (numFmtId 16 "d-mmm", centered).
$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = 10000
$ws.Cells.Item(9,3).Value = "Kim Hoang"
$ws.Cells.Item(9,4).NumberFormat = "d-mmm"
$ws.Cells.Item(9,4).HorizontalAlignment = -4108
$ws.Cells.Item(9,4).Value2 = 40433
$ws.Cells.Item(9,5).Value = "no submit"
$ws.Cells.Item(9,6).Value = "Waiting"

$ws.Cells.Item(10,1).Value = 9
$ws.Cells.Item(10,2).Value = 10000
$ws.Cells.Item(10,3).Value = "Thi Nguyen"
$ws.Cells.Item(10,4).NumberFormat = "d-mmm"
$ws.Cells.Item(10,4).HorizontalAlignment = -4108
$ws.Cells.Item(10,4).Value2 = 40433
$ws.Cells.Item(10,5).Value = "no submit"
$ws.Cells.Item(10,6).Value = "Waiting"

$ws.Cells.Item(11,1).Value = 10
$ws.Cells.Item(11,2).Value = 10000
$ws.Cells.Item(11,3).Value = "Kim Hoang"
$ws.Cells.Item(11,4).NumberFormat = "d-mmm"
$ws.Cells.Item(11,4).HorizontalAlignment = -4108
$ws.Cells.Item(11,4).Value2 = 40532
$ws.Cells.Item(11,5).Value = "late for meeting"
$ws.Cells.Item(11,6).Value = "Waiting"

$ws.Cells.Item(12,1).Value = 11
$ws.Cells.Item(12,2).Value = 10000
$ws.Cells.Item(12,3).Value = "Thi Nguyen"
$ws.Cells.Item(12,4).NumberFormat = "d-mmm"
$ws.Cells.Item(12,4).HorizontalAlignment = -4108
$ws.Cells.Item(12,4).Value2 = 40532
$ws.Cells.Item(12,5).Value = "rule 19"
$ws.Cells.Item(12,6).Value = "Waiting"

$ws.Cells.Item(13,1).Value = 12
$ws.Cells.Item(13,2).Value = 10000
$ws.Cells.Item(13,3).Value = "Tan Nguyen"
$ws.Cells.Item(13,4).NumberFormat = "d-mmm"
$ws.Cells.Item(13,4).HorizontalAlignment = -4108
$ws.Cells.Item(13,4).Value2 = 40532
$ws.Cells.Item(13,5).Value = "rule 19"
$ws.Cells.Item(13,6).Value = "Waiting"

$ws.Cells.Item(14,1).Value = 13
$ws.Cells.Item(14,2).Value = 10000
$ws.Cells.Item(14,3).Value = "Thang Le"
$ws.Cells.Item(14,4).NumberFormat = "d-mmm"
$ws.Cells.Item(14,4).HorizontalAlignment = -4108
$ws.Cells.Item(14,4).Value2 = 40532
$ws.Cells.Item(14,5).Value = "rule 19"
$ws.Cells.Item(14,6).Value = "Waiting"

# New totals row moved down to row 20, summing the now-larger data range
$ws.Cells.Item(20,2).Value = "sum"
$ws.Cells.Item(20,3).Formula = "=SUM(B2:B16)"

# Update selection to match the authored state
$ws.Range("G21").Select()

$wb.Save()
